# Add FEEDBACK MECHANISM to future scope
#
# Slide 8 ("Future Scope") has a bulleted text box ("TextBox 2") whose
# "Real-time integration ... with the sensors of acetabular reamer" bullet
# needs to be extended to also call out an efficient feed-back mechanism.
# The text box also needs to grow taller to accommodate the extra wrapped
# line of text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Locate the "TextBox 2" shape (Future Scope bullet list) by name so the
# script is resilient to shape ordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 2") {
        $shp = $s.Shapes.Item($i)
        break
    }
}

# Grow the text box's height (8687982 x 3462482 EMU -> 8687982 x 3831814 EMU)
# to fit the extra wrapped line introduced below. PowerPoint's Shape.Height
# is expressed in points (1 pt = 12700 EMU).
$shp.Height = 3831814 / 12700

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Find the paragraph that currently reads:
#   "Real-time integration with the sensors of acetabular reamer"
$targetParaIndex = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -like "*with the sensors of acetabular reamer*") {
        $targetParaIndex = $i
        break
    }
}

$para = $tr.Paragraphs($targetParaIndex)

# Find the run holding "with the sensors of acetabular reamer" within it.
$targetRun = $null
for ($i = 1; $i -le $para.Runs().Count; $i++) {
    if ($para.Runs($i).Text -like "*acetabular reamer*") {
        $targetRun = $para.Runs($i)
        break
    }
}

# Extend that run's text and append a new run carrying the new phrase.
$targetRun.Text = "with the sensors of acetabular reamer, with an efficient "
$newRun = $targetRun.InsertAfter("feed-back mechanism")
$newRun.Font.Size = 24
